$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply all cell content changes as plain text, forcing text interpretation
# so that values like "1.00" or "49.909.35" or "  -18.08%  " are preserved exactly
# (rather than being auto-coerced into numbers/percentages/dates by Excel).
function Set-TextValue($Cell, $Text) {
    $Cell.NumberFormat = "@"
    $Cell.Value = $Text
    $Cell.ClearFormats()
}

Set-TextValue $ws.Range('D2') '49.909.35'
Set-TextValue $ws.Range('E2') '  -18.08%  '
Set-TextValue $ws.Range('D3') '2.223.64'
Set-TextValue $ws.Range('E3') '  -23.60%  '
Set-TextValue $ws.Range('E4') '  +0.38%  '
Set-TextValue $ws.Range('D5') '413.42'
Set-TextValue $ws.Range('E5') '  -21.62%  '
Set-TextValue $ws.Range('D6') '112.44'
Set-TextValue $ws.Range('E6') '  -22.27%  '
Set-TextValue $ws.Range('D7') '1.00'
Set-TextValue $ws.Range('E7') '  +0.34%  '
Set-TextValue $ws.Range('D8') '0.442'
Set-TextValue $ws.Range('E8') '  -19.22%  '
Set-TextValue $ws.Range('D9') '2.224.78'
Set-TextValue $ws.Range('E9') '  -23.77%  '
Set-TextValue $ws.Range('D10') '4.97'
Set-TextValue $ws.Range('E10') '  -19.17%  '
Set-TextValue $ws.Range('E11') '  -22.79%  '
Set-TextValue $ws.Range('E12') '  -19.98%  '
Set-TextValue $ws.Range('E13') '  -7.98%  '
Set-TextValue $ws.Range('D14') '2.630.19'
Set-TextValue $ws.Range('E14') '  -23.02%  '
Set-TextValue $ws.Range('D15') '50.046.07'
Set-TextValue $ws.Range('E15') '  -17.76%  '
Set-TextValue $ws.Range('D16') '17.79'
Set-TextValue $ws.Range('E16') '  -21.11%  '
Set-TextValue $ws.Range('B17') 'WrappedEther'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D17') '2.262.40'
Set-TextValue $ws.Range('E17') '  -22.11%  '
Set-TextValue $ws.Range('B18') 'ShibaInu'
Set-TextValue $ws.Range('C18') 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Range('D18') '0.0000111'
Set-TextValue $ws.Range('E18') '  -21.75%  '
Set-TextValue $ws.Range('D19') '3.83'
Set-TextValue $ws.Range('E19') '  -21.64%  '
Set-TextValue $ws.Range('D20') '282.98'
Set-TextValue $ws.Range('E20') '  -20.19%  '
Set-TextValue $ws.Range('D21') '0.993'
Set-TextValue $ws.Range('E21') '  -0.86%  '
Set-TextValue $ws.Range('D22') '5.65'
Set-TextValue $ws.Range('E22') '  -1.30%  '
Set-TextValue $ws.Range('D23') '8.33'
Set-TextValue $ws.Range('E23') '  -28.26%  '
Set-TextValue $ws.Range('D24') '4.88'
Set-TextValue $ws.Range('E24') '  -25.15%  '
Set-TextValue $ws.Range('D25') '0.999'
Set-TextValue $ws.Range('E25') '  +0.01%  '
Set-TextValue $ws.Range('D26') '51.45'
Set-TextValue $ws.Range('E26') '  -20.78%  '
Set-TextValue $ws.Range('E27') '  -21.77%  '
Set-TextValue $ws.Range('D28') '2.333.94'
Set-TextValue $ws.Range('E28') '  -23.12%  '
Set-TextValue $ws.Range('B29') 'USDe'
Set-TextValue $ws.Range('C29') 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D29') '1.00'
Set-TextValue $ws.Range('E29') '  +0.10%  '
Set-TextValue $ws.Range('B30') 'Kaspa'
Set-TextValue $ws.Range('C30') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D30') '0.131'
Set-TextValue $ws.Range('E30') '  -26.28%  '
Set-TextValue $ws.Range('E31') '  -17.73%  '
Set-TextValue $ws.Range('D32') '0.0₃0620'
Set-TextValue $ws.Range('E32') '  -28.93%  '
Set-TextValue $ws.Range('D33') '140.37'
Set-TextValue $ws.Range('E33') '  -8.23%  '
Set-TextValue $ws.Range('E34') '  -18.43%  '
Set-TextValue $ws.Range('D35') '1.26'
Set-TextValue $ws.Range('E35') '  -25.41%  '
Set-TextValue $ws.Range('E36') '  -19.34%  '
Set-TextValue $ws.Range('D37') '1.00'
Set-TextValue $ws.Range('E37') '  +0.34%  '
Set-TextValue $ws.Range('D38') '3.19'
Set-TextValue $ws.Range('E38') '  -27.74%  '
Set-TextValue $ws.Range('D39') '31.41'
Set-TextValue $ws.Range('E39') '  -16.56%  '
Set-TextValue $ws.Range('B40') 'Fetch.AI'
Set-TextValue $ws.Range('C40') 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws.Range('D40') '0.724'
Set-TextValue $ws.Range('E40') '  -27.47%  '
Set-TextValue $ws.Range('B41') 'ImmutableX'
Set-TextValue $ws.Range('C41') 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D41') '0.932'
Set-TextValue $ws.Range('E41') '  -22.24%  '
Set-TextValue $ws.Range('E42') '  -2.05%  '
Set-TextValue $ws.Range('D43') '0.534'
Set-TextValue $ws.Range('E43') '  -18.16%  '
Set-TextValue $ws.Range('D44') '2.98'
Set-TextValue $ws.Range('E44') '  -19.67%  '
Set-TextValue $ws.Range('D45') '0.0476'
Set-TextValue $ws.Range('E45') '  -18.38%  '
Set-TextValue $ws.Range('D46') '1.823.88'
Set-TextValue $ws.Range('E46') '  -20.06%  '
Set-TextValue $ws.Range('E47') '  -27.21%  '
Set-TextValue $ws.Range('E48') '  -18.75%  '
Set-TextValue $ws.Range('D49') '0.0771'
Set-TextValue $ws.Range('E49') '  -16.11%  '
Set-TextValue $ws.Range('D50') '4.62'
Set-TextValue $ws.Range('E50') '  -5.90%  '
Set-TextValue $ws.Range('D51') '15.05'
Set-TextValue $ws.Range('E51') '  -26.11%  '
